$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "23.861.43"
$ws.Range("E2").Value = "  -2.93%  "
Set-TextValue $ws.Range("D3") "1.621.14"
$ws.Range("E3").Value = "  -3.20%  "
Set-TextValue $ws.Range("D4") "0.9998"
$ws.Range("E4").Value = "  -0.23%  "
Set-TextValue $ws.Range("D5") "308.69"
$ws.Range("E5").Value = "  -1.83%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.19%  "
Set-TextValue $ws.Range("D7") "0.3945"
$ws.Range("E7").Value = "  -0.14%  "
Set-TextValue $ws.Range("D8") "0.3840"
$ws.Range("E8").Value = "  -2.43%  "
$ws.Range("E9").Value = "  -0.17%  "
Set-TextValue $ws.Range("D10") "49.53"
$ws.Range("E10").Value = "  -2.39%  "
Set-TextValue $ws.Range("D11") "1.357"
$ws.Range("E11").Value = "  -2.60%  "
Set-TextValue $ws.Range("D12") "0.08456"
$ws.Range("E12").Value = "  -2.07%  "
Set-TextValue $ws.Range("D13") "23.72"
$ws.Range("E13").Value = "  -6.28%  "
Set-TextValue $ws.Range("D14") "7.052"
$ws.Range("E14").Value = "  -3.24%  "
Set-TextValue $ws.Range("D15") "7.584"
$ws.Range("E15").Value = "  -0.88%  "
Set-TextValue $ws.Range("D16") "0.00001281"
$ws.Range("E16").Value = "  -2.65%  "
Set-TextValue $ws.Range("D17") "1.618.98"
$ws.Range("E17").Value = "  -3.08%  "
Set-TextValue $ws.Range("D18") "93.90"
$ws.Range("E18").Value = "  +0.04%  "
Set-TextValue $ws.Range("D19") "0.06934"
$ws.Range("E19").Value = "  -1.20%  "
Set-TextValue $ws.Range("D20") "19.99"
$ws.Range("E20").Value = "  -5.96%  "
Set-TextValue $ws.Range("D21") "6.817"
$ws.Range("E21").Value = "  -3.56%  "
Set-TextValue $ws.Range("D22") "1.001"
$ws.Range("E22").Value = "  -0.18%  "
Set-TextValue $ws.Range("D23") "13.45"
$ws.Range("E23").Value = "  -3.06%  "
Set-TextValue $ws.Range("D24") "23.831.05"
$ws.Range("E24").Value = "  -3.02%  "
Set-TextValue $ws.Range("D25") "2.448"
$ws.Range("E25").Value = "  +4.30%  "
Set-TextValue $ws.Range("D26") "2.838"
$ws.Range("E26").Value = "  +2.31%  "
Set-TextValue $ws.Range("D27") "22.23"
$ws.Range("E27").Value = "  -3.19%  "
Set-TextValue $ws.Range("D28") "157.08"
$ws.Range("E28").Value = "  -1.30%  "
Set-TextValue $ws.Range("D29") "140.33"
$ws.Range("E29").Value = "  -3.50%  "
Set-TextValue $ws.Range("D30") "5.297"
$ws.Range("E30").Value = "  -9.14%  "
Set-TextValue $ws.Range("D31") "7.815"
$ws.Range("E31").Value = "  -5.00%  "
Set-TextValue $ws.Range("D32") "2.491"
$ws.Range("E32").Value = "  -2.31%  "
Set-TextValue $ws.Range("D33") "1.791.89"
$ws.Range("E33").Value = "  -3.61%  "
Set-TextValue $ws.Range("D34") "0.08127"
$ws.Range("E34").Value = "  -1.48%  "
Set-TextValue $ws.Range("D35") "0.9810"
$ws.Range("E35").Value = "  -1.49%  "
Set-TextValue $ws.Range("D36") "0.02875"
$ws.Range("E36").Value = "  -5.87%  "
Set-TextValue $ws.Range("D37") "6.588"
$ws.Range("E37").Value = "  -4.78%  "
Set-TextValue $ws.Range("D38") "0.2670"
$ws.Range("E38").Value = "  -4.64%  "
Set-TextValue $ws.Range("D39") "0.09146"
$ws.Range("E39").Value = "  -4.89%  "
Set-TextValue $ws.Range("D40") "10.38"
$ws.Range("E40").Value = "  +0.87%  "
Set-TextValue $ws.Range("D41") "13.61"
$ws.Range("E41").Value = "  +0.89%  "
Set-TextValue $ws.Range("D42") "1.428"
$ws.Range("E42").Value = "  -5.44%  "
Set-TextValue $ws.Range("D43") "0.7494"
$ws.Range("E43").Value = "  -4.74%  "
Set-TextValue $ws.Range("D44") "16.04"
$ws.Range("E44").Value = "  -2.41%  "
Set-TextValue $ws.Range("D45") "0.6924"
$ws.Range("E45").Value = "  -2.21%  "
Set-TextValue $ws.Range("D46") "2.474"
$ws.Range("E46").Value = "  -3.36%  "
Set-TextValue $ws.Range("D47") "4.078"
$ws.Range("E47").Value = "  -2.26%  "
Set-TextValue $ws.Range("D48") "1.0000"
$ws.Range("E48").Value = "  -0.26%  "
Set-TextValue $ws.Range("D49") "0.08238"
$ws.Range("E49").Value = "  -4.81%  "
Set-TextValue $ws.Range("D50") "134.71"
$ws.Range("E50").Value = "  -2.44%  "

# Row 51: Flow -> Tezos full row change
$ws.Range("B51").Value = "Tezos"
$ws.Range("C51").Value = "https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz"
Set-TextValue $ws.Range("D51") "1.387"
$ws.Range("E51").Value = "  +13.27%  "
